# Salary / materials sheet update:
#   - rows 2-8: replace the placeholder text with real content
#     (names, sum/dollar labels, prices, foreman), keeping numeric-looking
#     price strings ("123", "113", ...) stored as TEXT, like the source file.
#   - append a brand new row 9 with the same shape/style as the rows above.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: force a range's value to be stored as TEXT even when the string
# looks like a number (Excel would otherwise silently convert "113" -> 113).
# Tag the cell as text, write the value, then drop back to the default
# "Normal" style so no stray number-format style sticks to the cell.
function Set-TextValue($range, [string]$text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

# ---- Row 2 ----
$ws.Range("B2").Value = "dededede"
$ws.Range("C2").Value = "summ"
Set-TextValue $ws.Range("D2") "123"
$ws.Range("E2").Value = "effrfwwpkp"

# ---- Row 3 ----
$ws.Range("B3").Value = "dedede"
$ws.Range("C3").Value = "суммы"
Set-TextValue $ws.Range("D3") "113"
$ws.Range("E3").Value = "effrfwwpkp"

# ---- Row 4 ----
$ws.Range("B4").Value = "qwqeq"
$ws.Range("C4").Value = "/reload"
Set-TextValue $ws.Range("D4") "123"
$ws.Range("E4").Value = "effrfwwpkp"

# ---- Row 5 ----
$ws.Range("B5").Value = "Salary1"
$ws.Range("C5").Value = "суммы"
Set-TextValue $ws.Range("D5") "12"
$ws.Range("E5").Value = "effrfwwpkp"

# ---- Row 6 ----
$ws.Range("B6").Value = "ish xaqi1"
$ws.Range("C6").Value = "суммы"
Set-TextValue $ws.Range("D6") "87"
$ws.Range("E6").Value = "effrfwwpkp"

# ---- Row 7 ----
$ws.Range("B7").Value = "Hwjsjssj"
$ws.Range("C7").Value = "доллары"
Set-TextValue $ws.Range("D7") "123"
$ws.Range("E7").Value = "effrfwwpkp"

# ---- Row 8 ----
Set-TextValue $ws.Range("B8") "222222"
$ws.Range("C8").Value = "суммы"
Set-TextValue $ws.Range("D8") "100"
$ws.Range("E8").Value = "effrfwwpkp"

# ---- Row 9 (new) ----
$ws.Range("A9").Value = 7
$ws.Range("A8").Copy()
$ws.Range("A9").PasteSpecial(-4122)  # xlPasteFormats - match A2:A8 styling
$excel.CutCopyMode = $false

Set-TextValue $ws.Range("B9") "2222"
$ws.Range("C9").Value = "суммы"
Set-TextValue $ws.Range("D9") "12"
$ws.Range("E9").Value = "effrfwwpkp"
